$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5738271474838257
$ws.Range("B1").Value = 1.268808603286743
$ws.Range("C1").Value = 5.325573921203613
$ws.Range("D1").Value = 3.567200660705566
$ws.Range("E1").Value = 0.8273953795433044
